$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 180.3125  # H6
$ws.Cells.Item(6, 9).Value = 168.5  # I6
$ws.Cells.Item(6, 11).Value = 505.5  # K6
$ws.Cells.Item(6, 13).Value = -393.5  # M6
$ws.Cells.Item(17, 8).Value = 1853.75  # H17
$ws.Cells.Item(17, 10).Value = 2215  # J17
$ws.Cells.Item(17, 12).Value = 6645  # L17
$ws.Cells.Item(17, 14).Value = -6981  # N17
$ws.Cells.Item(70, 8).Value = 1600.375  # H70
$ws.Cells.Item(70, 9).Value = 1497.8889  # I70
$ws.Cells.Item(70, 10).Value = 1732.1428  # J70
$ws.Cells.Item(70, 11).Value = 4493.6667  # K70
$ws.Cells.Item(70, 12).Value = 5196.428400000001  # L70
$ws.Cells.Item(70, 13).Value = -4223.6667  # M70
$ws.Cells.Item(70, 14).Value = -5736.428400000001  # N70
$ws.Cells.Item(73, 8).Value = 1600.375  # H73
$ws.Cells.Item(73, 9).Value = 1497.8889  # I73
$ws.Cells.Item(73, 10).Value = 1732.1428  # J73
$ws.Cells.Item(73, 11).Value = 4493.6667  # K73
$ws.Cells.Item(73, 12).Value = 5196.428400000001  # L73
$ws.Cells.Item(73, 13).Value = -3557.6667  # M73
$ws.Cells.Item(73, 14).Value = -7068.428400000001  # N73
$ws.Cells.Item(74, 8).Value = 4279.6  # H74
$ws.Cells.Item(74, 9).Value = 2799.3333  # I74
$ws.Cells.Item(74, 10).Value = 6500  # J74
$ws.Cells.Item(74, 11).Value = 2799.3333  # K74
$ws.Cells.Item(74, 12).Value = 6500  # L74
$ws.Cells.Item(74, 13).Value = -1863.3333  # M74
$ws.Cells.Item(74, 14).Value = -8372  # N74
$ws.Cells.Item(77, 8).Value = 4279.6  # H77
$ws.Cells.Item(77, 9).Value = 2799.3333  # I77
$ws.Cells.Item(77, 10).Value = 6500  # J77
$ws.Cells.Item(77, 11).Value = 13996.6665  # K77
$ws.Cells.Item(77, 12).Value = 32500  # L77
$ws.Cells.Item(77, 13).Value = -9316.666499999999  # M77
$ws.Cells.Item(77, 14).Value = -41860  # N77
$ws.Cells.Item(80, 8).Value = 466.3684  # H80
$ws.Cells.Item(80, 10).Value = 500  # J80
$ws.Cells.Item(80, 12).Value = 1500  # L80
$ws.Cells.Item(80, 14).Value = -3496  # N80
$ws.Cells.Item(83, 8).Value = 466.3684  # H83
$ws.Cells.Item(83, 10).Value = 500  # J83
$ws.Cells.Item(83, 12).Value = 4500  # L83
$ws.Cells.Item(83, 14).Value = -14484  # N83
$ws.Cells.Item(111, 8).Value = 2869.5715  # H111
$ws.Cells.Item(111, 10).Value = 4149.25  # J111
$ws.Cells.Item(111, 12).Value = 12447.75  # L111
$ws.Cells.Item(111, 14).Value = -18581.75  # N111
$ws.Cells.Item(112, 8).Value = 2873.543  # H112
$ws.Cells.Item(112, 10).Value = 2895.2646  # J112
$ws.Cells.Item(112, 12).Value = 8685.793799999999  # L112
$ws.Cells.Item(112, 14).Value = -10901.7938  # N112
$ws.Cells.Item(135, 8).Value = 3307.5  # H135
$ws.Cells.Item(135, 9).Value = 3307.5  # I135
$ws.Cells.Item(135, 10).Value = 0  # J135
$ws.Cells.Item(135, 11).Value = 29767.5  # K135
$ws.Cells.Item(135, 12).Value = 0  # L135
$ws.Cells.Item(135, 13).Value = -27232.5  # M135
$ws.Cells.Item(135, 14).Value = $null  # N135
$ws.Cells.Item(141, 8).Value = 4622.643  # H141
$ws.Cells.Item(141, 9).Value = 4324.385  # I141
$ws.Cells.Item(141, 11).Value = 12973.155  # K141
$ws.Cells.Item(141, 13).Value = -7793.155000000001  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3159.42  # H32
$ws.Cells.Item(32, 9).Value = 2546.3257  # I32
$ws.Cells.Item(32, 10).Value = 6925.5713  # J32
$ws.Cells.Item(32, 11).Value = 2546.3257  # K32
$ws.Cells.Item(32, 12).Value = 6925.5713  # L32
$ws.Cells.Item(32, 13).Value = -2259.3257  # M32
$ws.Cells.Item(32, 14).Value = -7499.5713  # N32
$ws.Cells.Item(61, 8).Value = 3602.1428  # H61
$ws.Cells.Item(61, 9).Value = 2393.2173  # I61
$ws.Cells.Item(61, 10).Value = 5919.25  # J61
$ws.Cells.Item(61, 11).Value = 2393.2173  # K61
$ws.Cells.Item(61, 12).Value = 5919.25  # L61
$ws.Cells.Item(61, 13).Value = -2181.2173  # M61
$ws.Cells.Item(61, 14).Value = -6343.25  # N61
$ws.Cells.Item(101, 8).Value = 160068  # H101
$ws.Cells.Item(101, 10).Value = 160068  # J101
$ws.Cells.Item(101, 12).Value = 160068  # L101
$ws.Cells.Item(101, 14).Value = -166558  # N101
$ws.Cells.Item(133, 8).Value = 60000  # H133
$ws.Cells.Item(133, 10).Value = 60000  # J133
$ws.Cells.Item(133, 12).Value = 60000  # L133
$ws.Cells.Item(133, 14).Value = -65060  # N133
$ws.Cells.Item(136, 8).Value = 3602.1428  # H136
$ws.Cells.Item(136, 9).Value = 2393.2173  # I136
$ws.Cells.Item(136, 10).Value = 5919.25  # J136
$ws.Cells.Item(136, 11).Value = 7179.651899999999  # K136
$ws.Cells.Item(136, 12).Value = 17757.75  # L136
$ws.Cells.Item(136, 13).Value = -4629.651899999999  # M136
$ws.Cells.Item(136, 14).Value = -22857.75  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 66669140  # H20
$ws.Cells.Item(20, 9).Value = 90911360  # I20
$ws.Cells.Item(20, 10).Value = 3024.75  # J20
$ws.Cells.Item(20, 11).Value = 90911360  # K20
$ws.Cells.Item(20, 12).Value = 3024.75  # L20
$ws.Cells.Item(20, 13).Value = -90911113  # M20
$ws.Cells.Item(20, 14).Value = -3518.75  # N20
$ws.Cells.Item(86, 8).Value = 2695  # H86
$ws.Cells.Item(86, 9).Value = 2543.6428  # I86
$ws.Cells.Item(86, 10).Value = 3118.8  # J86
$ws.Cells.Item(86, 11).Value = 2543.6428  # K86
$ws.Cells.Item(86, 12).Value = 3118.8  # L86
$ws.Cells.Item(86, 13).Value = -1420.6428  # M86
$ws.Cells.Item(86, 14).Value = -5364.8  # N86
$ws.Cells.Item(89, 8).Value = 2695  # H89
$ws.Cells.Item(89, 9).Value = 2543.6428  # I89
$ws.Cells.Item(89, 10).Value = 3118.8  # J89
$ws.Cells.Item(89, 11).Value = 12718.214  # K89
$ws.Cells.Item(89, 12).Value = 15594  # L89
$ws.Cells.Item(89, 13).Value = -7102.214  # M89
$ws.Cells.Item(89, 14).Value = -26826  # N89
$ws.Cells.Item(134, 8).Value = 13514678  # H134
$ws.Cells.Item(134, 9).Value = 16667920  # I134
$ws.Cells.Item(134, 11).Value = 50003760  # K134
$ws.Cells.Item(134, 13).Value = -50001225  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 884  # H22
$ws.Cells.Item(22, 10).Value = 1621.5555  # J22
$ws.Cells.Item(22, 12).Value = 1621.5555  # L22
$ws.Cells.Item(22, 14).Value = -2321.5555  # N22
$ws.Cells.Item(31, 8).Value = 2218.4358  # H31
$ws.Cells.Item(31, 9).Value = 2066.5925  # I31
$ws.Cells.Item(31, 10).Value = 2560.0833  # J31
$ws.Cells.Item(31, 11).Value = 2066.5925  # K31
$ws.Cells.Item(31, 12).Value = 2560.0833  # L31
$ws.Cells.Item(31, 13).Value = -1771.5925  # M31
$ws.Cells.Item(31, 14).Value = -3150.0833  # N31
$ws.Cells.Item(34, 8).Value = 2218.4358  # H34
$ws.Cells.Item(34, 9).Value = 2066.5925  # I34
$ws.Cells.Item(34, 10).Value = 2560.0833  # J34
$ws.Cells.Item(34, 11).Value = 2066.5925  # K34
$ws.Cells.Item(34, 12).Value = 2560.0833  # L34
$ws.Cells.Item(34, 13).Value = -1864.5925  # M34
$ws.Cells.Item(34, 14).Value = -2964.0833  # N34
$ws.Cells.Item(134, 8).Value = 906.82355  # H134
$ws.Cells.Item(134, 9).Value = 906.82355  # I134
$ws.Cells.Item(134, 11).Value = 2720.47065  # K134
$ws.Cells.Item(134, 13).Value = -185.4706499999998  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 497  # H7
$ws.Cells.Item(7, 9).Value = 0  # I7
$ws.Cells.Item(7, 11).Value = 0  # K7
$ws.Cells.Item(7, 13).Value = $null  # M7
$ws.Cells.Item(45, 8).Value = 3033  # H45
$ws.Cells.Item(45, 10).Value = 3033  # J45
$ws.Cells.Item(45, 12).Value = 9099  # L45
$ws.Cells.Item(45, 14).Value = -10163  # N45
$ws.Cells.Item(92, 8).Value = 299  # H92
$ws.Cells.Item(92, 9).Value = 299  # I92
$ws.Cells.Item(92, 11).Value = 897  # K92
$ws.Cells.Item(92, 13).Value = 351  # M92
$ws.Cells.Item(131, 8).Value = 1598.6  # H131
$ws.Cells.Item(131, 10).Value = 2915  # J131
$ws.Cells.Item(131, 12).Value = 8745  # L131
$ws.Cells.Item(131, 14).Value = -18825  # N131
$ws.Cells.Item(132, 8).Value = 1314.44  # H132
$ws.Cells.Item(132, 9).Value = 1251.1875  # I132
$ws.Cells.Item(132, 10).Value = 1426.8889  # J132
$ws.Cells.Item(132, 11).Value = 11260.6875  # K132
$ws.Cells.Item(132, 12).Value = 12842.0001  # L132
$ws.Cells.Item(132, 13).Value = -8730.6875  # M132
$ws.Cells.Item(132, 14).Value = -17902.0001  # N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(23, 8).Value = 0  # H23
$ws.Cells.Item(23, 9).Value = 0  # I23
$ws.Cells.Item(23, 11).Value = 0  # K23
$ws.Cells.Item(23, 13).Value = $null  # M23
$ws.Cells.Item(80, 8).Value = 3912.074  # H80
$ws.Cells.Item(80, 9).Value = 3851  # I80
$ws.Cells.Item(80, 10).Value = 3960.9333  # J80
$ws.Cells.Item(80, 11).Value = 3851  # K80
$ws.Cells.Item(80, 12).Value = 3960.9333  # L80
$ws.Cells.Item(80, 13).Value = -2853  # M80
$ws.Cells.Item(80, 14).Value = -5956.933300000001  # N80
$ws.Cells.Item(83, 8).Value = 3912.074  # H83
$ws.Cells.Item(83, 9).Value = 3851  # I83
$ws.Cells.Item(83, 10).Value = 3960.9333  # J83
$ws.Cells.Item(83, 11).Value = 19255  # K83
$ws.Cells.Item(83, 12).Value = 19804.6665  # L83
$ws.Cells.Item(83, 13).Value = -14263  # M83
$ws.Cells.Item(83, 14).Value = -29788.6665  # N83
$ws.Cells.Item(136, 8).Value = 71110.734  # H136
$ws.Cells.Item(136, 10).Value = 71110.734  # J136
$ws.Cells.Item(136, 12).Value = 213332.202  # L136
$ws.Cells.Item(136, 14).Value = -218432.202  # N136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1483  # H46
$ws.Cells.Item(46, 9).Value = 1162  # I46
$ws.Cells.Item(46, 11).Value = 1162  # K46
$ws.Cells.Item(46, 13).Value = -974  # M46
$ws.Cells.Item(55, 8).Value = 1222.7084  # H55
$ws.Cells.Item(55, 9).Value = 517.25  # I55
$ws.Cells.Item(55, 11).Value = 517.25  # K55
$ws.Cells.Item(55, 13).Value = -344.25  # M55
$ws.Cells.Item(68, 8).Value = 3624.5  # H68
$ws.Cells.Item(68, 9).Value = 2405.2  # I68
$ws.Cells.Item(68, 10).Value = 9721  # J68
$ws.Cells.Item(68, 11).Value = 2405.2  # K68
$ws.Cells.Item(68, 12).Value = 9721  # L68
$ws.Cells.Item(68, 13).Value = -1656.2  # M68
$ws.Cells.Item(68, 14).Value = -11219  # N68
$ws.Cells.Item(71, 8).Value = 3624.5  # H71
$ws.Cells.Item(71, 9).Value = 2405.2  # I71
$ws.Cells.Item(71, 10).Value = 9721  # J71
$ws.Cells.Item(71, 11).Value = 12026  # K71
$ws.Cells.Item(71, 12).Value = 48605  # L71
$ws.Cells.Item(71, 13).Value = -8282  # M71
$ws.Cells.Item(71, 14).Value = -56093  # N71
$ws.Cells.Item(133, 8).Value = 66000  # H133
$ws.Cells.Item(133, 10).Value = 66000  # J133
$ws.Cells.Item(133, 12).Value = 66000  # L133
$ws.Cells.Item(133, 14).Value = -71060  # N133
$ws.Cells.Item(136, 8).Value = 5588.3257  # H136
$ws.Cells.Item(136, 9).Value = 4353.8184  # I136
$ws.Cells.Item(136, 11).Value = 13061.4552  # K136
$ws.Cells.Item(136, 13).Value = -10511.4552  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 4993.3335  # H8
$ws.Cells.Item(8, 9).Value = 4993.3335  # I8
$ws.Cells.Item(8, 10).Value = 0  # J8
$ws.Cells.Item(8, 11).Value = 4993.3335  # K8
$ws.Cells.Item(8, 12).Value = 0  # L8
$ws.Cells.Item(8, 13).Value = -4853.3335  # M8
$ws.Cells.Item(8, 14).Value = $null  # N8
$ws.Cells.Item(11, 8).Value = 74998.336  # H11
$ws.Cells.Item(11, 9).Value = 74998.336  # I11
$ws.Cells.Item(11, 11).Value = 74998.336  # K11
$ws.Cells.Item(11, 13).Value = -74856.336  # M11
$ws.Cells.Item(96, 8).Value = 7409214.5  # H96
$ws.Cells.Item(96, 9).Value = 12347253  # I96
$ws.Cells.Item(96, 10).Value = 2157.5  # J96
$ws.Cells.Item(96, 11).Value = 12347253  # K96
$ws.Cells.Item(96, 12).Value = 2157.5  # L96
$ws.Cells.Item(96, 13).Value = -12345880  # M96
$ws.Cells.Item(96, 14).Value = -4903.5  # N96
$ws.Cells.Item(108, 8).Value = 69990  # H108
$ws.Cells.Item(108, 10).Value = 69990  # J108
$ws.Cells.Item(108, 12).Value = 69990  # L108
$ws.Cells.Item(108, 14).Value = -77670  # N108
$ws.Cells.Item(138, 8).Value = 59994.5  # H138
$ws.Cells.Item(138, 10).Value = 59994.5  # J138
$ws.Cells.Item(138, 12).Value = 59994.5  # L138
$ws.Cells.Item(138, 14).Value = -70274.5  # N138
$ws.Cells.Item(139, 8).Value = 100000  # H139
$ws.Cells.Item(139, 10).Value = 100000  # J139
$ws.Cells.Item(139, 12).Value = 100000  # L139
$ws.Cells.Item(139, 14).Value = -110280  # N139

